$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - reorder/rewrite
$ws.Range("A1").Value = "D"
$ws.Range("B1").Value = "t-1"
$ws.Range("C1").Value = "t-2"
$ws.Range("D1").Value = "T"
$ws.Range("E1").Value = "X"
$ws.Range("F1").Value = "Y"

# Clear old header text (keep styling) for cells no longer holding labels
$ws.Range("G1").ClearContents()
$ws.Range("H1").ClearContents()
$ws.Range("J1").ClearContents()
$ws.Range("K1").ClearContents()
$ws.Range("L1").ClearContents()
$ws.Range("M1").ClearContents()

# Clear row 2 extra calculation cells
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()

# Update selection to match the new active range
$ws.Range("G1:M1").Select()
